# The commit replaces `<w:t xml:space="preserve"/>` with `<w:t/>` for every
# *empty* text run that sits right after a loop-tag run in the 4x3 table
# (12 runs total: 4 rows x 3 columns x 1 empty-preserve run per cell).
# The Word OM has no "Run" object and empty runs occupy zero characters, so
# they cannot be targeted with Range(start,end)/Find. Instead we rebuild each
# cell paragraph's exact OOXML (all runs/attrs unchanged) via InsertXML, only
# dropping xml:space="preserve" from the runs whose <w:t> is empty.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$cellXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009C697E" w:rsidRPr="00994758" w:rsidRDefault="00306398" w:rsidP="00EF403C"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00F5114B"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="008A103F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">John</w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="008A103F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tbl.Cell(1,1).Range.Paragraphs.Item(1).Range.InsertXML($cellXml)

$cellXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009C697E" w:rsidRPr="00994758" w:rsidRDefault="0055408B" w:rsidP="00EF403C"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Doe</w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00B7631C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tbl.Cell(1,2).Range.Paragraphs.Item(1).Range.InsertXML($cellXml)

$cellXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009C697E" w:rsidRPr="00994758" w:rsidRDefault="008A103F" w:rsidP="00EF403C"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">+33647874513</w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="009C697E" w:rsidRPr="00994758"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="009C697E" w:rsidRPr="00994758"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tbl.Cell(1,3).Range.Paragraphs.Item(1).Range.InsertXML($cellXml)

$cellXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009C697E" w:rsidRPr="00994758" w:rsidRDefault="00306398" w:rsidP="00EF403C"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00F5114B"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="008A103F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Jane</w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="008A103F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tbl.Cell(2,1).Range.Paragraphs.Item(1).Range.InsertXML($cellXml)

$cellXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009C697E" w:rsidRPr="00994758" w:rsidRDefault="0055408B" w:rsidP="00EF403C"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Doe</w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00B7631C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tbl.Cell(2,2).Range.Paragraphs.Item(1).Range.InsertXML($cellXml)

$cellXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009C697E" w:rsidRPr="00994758" w:rsidRDefault="008A103F" w:rsidP="00EF403C"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">+33454540124</w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="009C697E" w:rsidRPr="00994758"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="009C697E" w:rsidRPr="00994758"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tbl.Cell(2,3).Range.Paragraphs.Item(1).Range.InsertXML($cellXml)

$cellXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009C697E" w:rsidRPr="00994758" w:rsidRDefault="00306398" w:rsidP="00EF403C"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00F5114B"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="008A103F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Phil</w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="008A103F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tbl.Cell(3,1).Range.Paragraphs.Item(1).Range.InsertXML($cellXml)

$cellXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009C697E" w:rsidRPr="00994758" w:rsidRDefault="0055408B" w:rsidP="00EF403C"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Kiel</w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00B7631C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tbl.Cell(3,2).Range.Paragraphs.Item(1).Range.InsertXML($cellXml)

$cellXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009C697E" w:rsidRPr="00994758" w:rsidRDefault="008A103F" w:rsidP="00EF403C"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">+44578451245</w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="009C697E" w:rsidRPr="00994758"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="009C697E" w:rsidRPr="00994758"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tbl.Cell(3,3).Range.Paragraphs.Item(1).Range.InsertXML($cellXml)

$cellXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009C697E" w:rsidRPr="00994758" w:rsidRDefault="00306398" w:rsidP="00EF403C"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00F5114B"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="008A103F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Dave</w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="008A103F"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tbl.Cell(4,1).Range.Paragraphs.Item(1).Range.InsertXML($cellXml)

$cellXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009C697E" w:rsidRPr="00994758" w:rsidRDefault="0055408B" w:rsidP="00EF403C"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Sto</w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00B7631C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tbl.Cell(4,2).Range.Paragraphs.Item(1).Range.InsertXML($cellXml)

$cellXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="009C697E" w:rsidRPr="00994758" w:rsidRDefault="008A103F" w:rsidP="00EF403C"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">+44548787984</w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="009C697E" w:rsidRPr="00994758"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:r w:rsidR="00EF403C"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="009C697E" w:rsidRPr="00994758"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t></w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tbl.Cell(4,3).Range.Paragraphs.Item(1).Range.InsertXML($cellXml)

